# Normalizes the "Recorded By" (column G) audit-trail strings on the
# "Session Analysis Results" sheet: the comma-separated list of recorders
# is reversed, except that "backup@backdoor.com" always keeps its original
# position in the list (it is pinned; only the remaining entries around it
# are reversed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pinned = "backup@backdoor.com"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $raw = $cell.Value2

    if ($raw -eq $null) { continue }
    if ($raw -eq "") { continue }

    $parts = $raw -split ", "
    if ($parts.Count -le 1) { continue }

    $pinnedIndex = -1
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($parts[$i] -eq $pinned) {
            $pinnedIndex = $i
        }
    }

    if ($pinnedIndex -ge 0) {
        $rest = @()
        for ($i = 0; $i -lt $parts.Count; $i++) {
            if ($i -ne $pinnedIndex) { $rest += $parts[$i] }
        }
        $restRev = @()
        for ($i = $rest.Count - 1; $i -ge 0; $i--) { $restRev += $rest[$i] }

        $newParts = @()
        for ($i = 0; $i -lt $pinnedIndex; $i++) { $newParts += $restRev[$i] }
        $newParts += $pinned
        for ($i = $pinnedIndex; $i -lt $restRev.Count; $i++) { $newParts += $restRev[$i] }
    } else {
        $newParts = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) { $newParts += $parts[$i] }
    }

    $newValue = $newParts -join ", "
    $cell.Value = $newValue
}
